$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("historico")

$data = @(
    ,@(197, '09/01/2026 10:51', 439, 'Conhecimentos Específicos', 'Logística', 0)
    ,@(198, '09/01/2026 10:55', 443, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(199, '09/01/2026 10:57', 1810, 'Conhecimentos Específicos', 'Logística', 0)
    ,@(200, '09/01/2026 10:59', 828, 'Conhecimentos Específicos', 'Logística', 0)
    ,@(201, '09/01/2026 11:02', 757, 'Conhecimentos Específicos', 'Logística', 0)
    ,@(202, '09/01/2026 11:06', 1811, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(203, '09/01/2026 11:08', 814, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(204, '09/01/2026 11:08', 434, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(205, '09/01/2026 11:09', 1808, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(206, '09/01/2026 11:09', 438, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(207, '09/01/2026 11:10', 1790, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(208, '09/01/2026 11:13', 1806, 'Conhecimentos Específicos', 'Logística', 0)
    ,@(209, '09/01/2026 11:13', 853, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(210, '09/01/2026 11:15', 1791, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(211, '09/01/2026 11:16', 444, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(212, '09/01/2026 11:17', 435, 'Conhecimentos Específicos', 'Logística', 0)
    ,@(213, '09/01/2026 11:18', 780, 'Conhecimentos Específicos', 'Logística', 0)
    ,@(214, '09/01/2026 11:19', 775, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(215, '09/01/2026 11:20', 759, 'Conhecimentos Específicos', 'Logística', 1)
    ,@(216, '09/01/2026 11:21', 1804, 'Conhecimentos Específicos', 'Logística', 1)
)

$startRow = 184
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

Write-Output "Added $($data.Count) rows"
